# Reattempt "Is Anagram" and "Best Time to Buy and Sell Stock" (and the
# knock-on star-count bumps on a few neighbouring rows), per the commit:
# "reattemp Best time to buy & sell and valid anagram"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - Is Anagram: 1* -> 3*
$ws.Range("E12").Value = "3*"

# Row 14 - Best Time to Buy and Sell Stock: 4* -> 5*
$ws.Range("E14").Value = "5*"

# Row 20 - Merge Two Sorted Lists: 5* -> 6*
$ws.Range("E20").Value = "6*"

# Row 21 - Reverse Linked List (Recursively): 3* -> 4*
$ws.Range("E21").Value = "4*"

# Row 22 - Has Cycle: 5* -> 6*
$ws.Range("E22").Value = "6*"

# Row heights 19.5 -> 18.75 for rows 1, 2, 5, 6, 7 (matches the autofit
# recalculation that accompanied the edits above)
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 18.75
$ws.Rows.Item(5).RowHeight = 18.75
$ws.Rows.Item(6).RowHeight = 18.75
$ws.Rows.Item(7).RowHeight = 18.75
